$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-22 22:18:21'
$ws.Range("H2").Value = '''41%'
$ws.Range("O2").Value = '6.0 °C'
$ws.Range("E3").Value = '2026-02-22 22:18:24'
$ws.Range("O3").Value = '4.3 °C'
$ws.Range("E4").Value = '2026-02-22 22:18:26'
$ws.Range("O4").Value = '12.2 °C'
$ws.Range("E5").Value = '2026-02-22 22:18:28'
$ws.Range("O5").Value = '5.9 °C'
$ws.Range("E6").Value = '2026-02-22 22:18:31'
$ws.Range("O6").Value = '13.0 °C'
$ws.Range("E7").Value = '2026-02-22 22:18:34'
$ws.Range("H7").Value = '''60%'
$ws.Range("J7").Value = '1027.3 hPa'
$ws.Range("E8").Value = '2026-02-22 22:18:36'
$ws.Range("H8").Value = '''48%'
$ws.Range("E9").Value = '2026-02-22 22:18:39'
$ws.Range("E10").Value = '2026-02-22 22:18:41'
$ws.Range("E11").Value = '2026-02-22 22:18:42'
$ws.Range("O11").Value = '8.6 °C'
$ws.Range("E12").Value = '2026-02-22 22:18:43'
$ws.Range("E13").Value = '2026-02-22 22:18:44'
$ws.Range("E14").Value = '2026-02-22 22:18:45'
$ws.Range("E15").Value = '2026-02-22 22:18:46'
$ws.Range("E16").Value = '2026-02-22 22:18:48'
$ws.Range("L16").Value = '20.5 km/h - 236º 21:59 TU'
$ws.Range("E17").Value = '2026-02-22 22:18:49'
$ws.Range("N17").Value = '7.6 °C 21:46 TU'
$ws.Range("O17").Value = '9.9 °C'
$ws.Range("E18").Value = '2026-02-22 22:18:50'
$ws.Range("O18").Value = '10.0 °C'
$ws.Range("E19").Value = '2026-02-22 22:18:51'
$ws.Range("H19").Value = '''48%'
$ws.Range("O19").Value = '12.0 °C'
$ws.Range("E20").Value = '2026-02-22 22:18:52'
$ws.Range("E21").Value = '2026-02-22 22:18:53'
$ws.Range("E22").Value = '2026-02-22 22:18:56'
$ws.Range("E23").Value = '2026-02-22 22:18:58'
$ws.Range("E24").Value = '2026-02-22 22:19:01'
$ws.Range("E25").Value = '2026-02-22 22:19:03'
$ws.Range("K25").Value = '16.6 MJ/m2'
$ws.Range("E26").Value = '2026-02-22 22:19:05'
$ws.Range("H26").Value = '''38%'
$ws.Range("J26").Value = '1026.1 hPa'
$ws.Range("E27").Value = '2026-02-22 22:19:08'
$ws.Range("H27").Value = '''27%'
$ws.Range("E28").Value = '2026-02-22 22:19:10'
$ws.Range("H28").Value = '''65%'
$ws.Range("O28").Value = '10.3 °C'
$ws.Range("E29").Value = '2026-02-22 22:19:13'
$ws.Range("K29").Value = '15.3 MJ/m2'
$ws.Range("E30").Value = '2026-02-22 22:19:15'
$ws.Range("H30").Value = '''74%'
$ws.Range("O30").Value = '12.1 °C'
$ws.Range("E31").Value = '2026-02-22 22:19:18'
$ws.Range("K31").Value = '15.1 MJ/m2'
$ws.Range("L31").Value = '22.0 km/h - 11º 21:37 TU'
$ws.Range("E32").Value = '2026-02-22 22:19:20'
$ws.Range("H32").Value = '''72%'
$ws.Range("O32").Value = '5.7 °C'
$ws.Range("E33").Value = '2026-02-22 22:19:22'
$ws.Range("H33").Value = '''49%'
$ws.Range("E34").Value = '2026-02-22 22:19:25'
$ws.Range("E35").Value = '2026-02-22 22:19:28'
$ws.Range("J35").Value = '1028.4 hPa'
$ws.Range("K35").Value = '16.2 MJ/m2'
$ws.Range("E36").Value = '2026-02-22 22:19:30'
$ws.Range("H36").Value = '''77%'
$ws.Range("O36").Value = '11.6 °C'
$ws.Range("E37").Value = '2026-02-22 22:19:33'
$ws.Range("J37").Value = '1029.8 hPa'
$ws.Range("E38").Value = '2026-02-22 22:19:35'
$ws.Range("E39").Value = '2026-02-22 22:19:38'
$ws.Range("E40").Value = '2026-02-22 22:19:40'
$ws.Range("O40").Value = '9.8 °C'
$ws.Range("E41").Value = '2026-02-22 22:19:43'
$ws.Range("E42").Value = '2026-02-22 22:19:45'
$ws.Range("H42").Value = '''82%'
$ws.Range("O42").Value = '10.5 °C'
$ws.Range("E43").Value = '2026-02-22 22:19:47'
$ws.Range("E44").Value = '2026-02-22 22:19:49'
$ws.Range("E45").Value = '2026-02-22 22:19:52'
$ws.Range("J45").Value = '1029.1 hPa'
$ws.Range("O45").Value = '8.6 °C'
$ws.Range("E46").Value = '2026-02-22 22:19:55'
$ws.Range("H46").Value = '''76%'
$ws.Range("O46").Value = '8.9 °C'
